$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = -13.15399999999999
$ws.Range("A9").Value = -21.97940000000002
$ws.Range("C12").Value = -11.40789999999999
$ws.Range("E13").Value = 16.6507
$ws.Range("D15").Value = -8.605599999999994
$ws.Range("E16").Value = 16.4519
$ws.Range("A18").Value = -22.08860000000001
$ws.Range("A20").Value = -20.5933
$ws.Range("E20").Value = 15.89169999999999
$ws.Range("E24").Value = 16.71300000000001
$ws.Range("C26").Value = -12.7188
$ws.Range("A27").Value = -21.78439999999999
$ws.Range("C27").Value = -12.8195
$ws.Range("C29").Value = -11.20320000000001
$ws.Range("C37").Value = -13.83669999999999
$ws.Range("C38").Value = -13.2552
$ws.Range("D38").Value = -8.898199999999989
$ws.Range("E39").Value = 16.0125
$ws.Range("D44").Value = -7.317400000000003
$ws.Range("E48").Value = 17.30990000000001
$ws.Range("C51").Value = -12.0709
$ws.Range("D51").Value = -7.1201
$ws.Range("E52").Value = 17.1387
$ws.Range("C55").Value = -13.57950000000001
$ws.Range("E56").Value = 16.56000000000001
$ws.Range("D57").Value = -8.334400000000002
$ws.Range("D63").Value = -7.709599999999996
$ws.Range("A69").Value = -21.97340000000001
$ws.Range("C69").Value = -11.7362
$ws.Range("C70").Value = -11.9592
$ws.Range("D70").Value = -7.845799999999998
$ws.Range("A76").Value = -20.04869999999999
$ws.Range("A82").Value = -22.15590000000001
$ws.Range("C83").Value = -13.9588
$ws.Range("E84").Value = 16.7305
$ws.Range("D99").Value = -8.3066
$ws.Range("E100").Value = 16.4444
$ws.Range("E101").Value = 16.82440000000001
$ws.Range("C102").Value = -14.0645
